# Sprint 11 alignment: studyDesignProcedures gains two new columns
# ("procedureName" and "procedureDescription") inserted right after the
# existing "procedureType" column, and the user ends up with that sheet
# as the active tab/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studyDesignProcedures")

# Insert two new blank columns before the current column C
# (procedureCode), pushing procedureCode/procedureIsConditional/
# procedureIsConditionalReason two columns to the right. The new
# columns inherit column B's formatting (same as a normal Excel
# "Insert Column" action run from a selection spanning two columns).
$ws.Range("C1:D1").EntireColumn.Insert() | Out-Null

# Column B's width carries over automatically to the freshly inserted
# column C; give column C that same width explicitly and give the new
# column D its own (slightly wider) custom width.
$ws.Columns.Item(3).ColumnWidth = 16.67
$ws.Columns.Item(4).ColumnWidth = 20.67

# New header row 1 labels (leading ' keeps these plain, left/top
# aligned labels matching the style already used by the neighboring
# header cells).
$ws.Range("C1").Value = "'procedureName"
$ws.Range("D1").Value = "'procedureDescription"

# New sample data rows.
$ws.Range("C2").Value = "Test8"
$ws.Range("D2").Value = "Test Eight"
$ws.Range("C3").Value = "Test9"
$ws.Range("D3").Value = "Test Nine"

# The edit session ends with studyDesignProcedures as the active sheet
# (previously studyDesign was active) and a selection on D9.
$ws.Activate() | Out-Null
$ws.Range("D9").Select() | Out-Null
